# Update the "enhance easyocr text testing results" changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update EasyOCR row: Error Rate on Alphabets (D4) value
$ws.Range("D4").Value = 0.043

# Update EasyOCR row: Misinterpret Alphabets (F4) text
$ws.Range("F4").Value = "misinterpret l to i`nmisinterpret h to n`nmisinterpret f to t`nmisinterpret d to a`nmiss y in the end`nmiss v"

# Update the selected cell in the sheet view
$ws.Range("D7").Select() | Out-Null
